$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# C2: "LAB DIS - 1" -> "1-32"
$ws.Range("C2").Value = "1-32"

# E2: "LAB LTO / PALESTRA 1 / PALESTRA 2 / SALA PESI" -> "LTO / PAL1 / PAL2 / PALF"
$ws.Range("E2").Value = "LTO" + $nl + "PAL1" + $nl + "PAL2" + $nl + "PALF"

# F2: replace trailing "BASKIN" with "PBAS"
$ws.Range("F2").Value = "1-14" + $nl + "1-15" + $nl + "1-17" + $nl + "1-90" + $nl + "1-91" + $nl + "1-92" + $nl + "1-93" + $nl + "1-94" + $nl + "1-95" + $nl + "1-96" + $nl + "PBAS"

# G3: insert E0-1/E0-2/E0-3, M0-1/M0-2/M0-3, T0-1, I0-1..I0-5 entries into the list
$ws.Range("G3").Value = "E0-1" + $nl + "E0-2" + $nl + "E0-3" + $nl + "E1-7" + $nl + "E1-8" + $nl + "M0-1" + $nl + "M0-2" + $nl + "M0-3" + $nl + "M1-8" + $nl + "M1-11" + $nl + "M1-18" + $nl + "T0-1" + $nl + "T1-6" + $nl + "T1-7" + $nl + "T1-12" + $nl + "T1-13" + $nl + "T1-14" + $nl + "I0-1" + $nl + "I0-2" + $nl + "I0-3" + $nl + "I0-4" + $nl + "I0-5" + $nl + "I1-1" + $nl + "I1-2" + $nl + "I1-3" + $nl + "I1-6" + $nl + "I1-13"

# B4: replace trailing "LAB DIS - 2" with "2-8"
$ws.Range("B4").Value = "2-1" + $nl + "2-2" + $nl + "2-5" + $nl + "2-6" + $nl + "2-7" + $nl + "2-8"
